$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel does not
# auto-convert them to numbers (losing trailing zeros / exact text).
$forceTextCells = @("D5", "D6", "D13", "D14", "D21", "D25", "D26", "D27", "D29", "D31", "D35", "D38", "D47", "D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = '66.497.90'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '2.579.97'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("D5").Value = '582.41'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = '165.83'
$ws.Range("E6").Value = '  +1.05%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").Value = '2.578.34'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  -3.41%  '
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '5.18'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '26.80'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").Value = '3.050.64'
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").Value = '66.354.68'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").Value = '2.567.93'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("E19").Value = '  -4.52%  '
$ws.Range("E20").Value = '  -3.31%  '
$ws.Range("D21").Value = '351.78'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = '1.89'
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '10.03'
$ws.Range("E26").Value = '  -7.82%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '68.97'
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").Value = '2.709.86'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '0.0₃0988'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("D31").Value = '536.86'
$ws.Range("E31").Value = '  -2.23%  '
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  -2.94%  '
$ws.Range("D38").Value = '156.77'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("D47").Value = '149.04'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '0.0759'
$ws.Range("E51").Value = '  -1.58%  '
